$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.23005485534668
$ws.Range("B1").Value = 2.532801389694214
$ws.Range("C1").Value = 4.949195861816406
$ws.Range("D1").Value = 2.261281490325928
$ws.Range("E1").Value = 1.062237024307251
